$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 382.66666
$ws.Range("I42").Value = 382.66666
$ws.Range("K42").Value = 1147.99998
$ws.Range("M42").Value = -917.9999800000001

# Row 53
$ws.Range("H53").Value = 429.14285
$ws.Range("I53").Value = 95.833336
$ws.Range("K53").Value = 95.833336
$ws.Range("M53").Value = 541.166664

# Row 61
$ws.Range("H61").Value = 345.6
$ws.Range("I61").Value = 344.75
$ws.Range("J61").Value = 349
$ws.Range("K61").Value = 1034.25
$ws.Range("L61").Value = 1047
$ws.Range("M61").Value = -862.25
$ws.Range("N61").Value = -1391

# Row 70
$ws.Range("H70").Value = 6105.8696
$ws.Range("I70").Value = 1481.6
$ws.Range("K70").Value = 4444.799999999999
$ws.Range("M70").Value = -4174.799999999999

# Row 73
$ws.Range("H73").Value = 6105.8696
$ws.Range("I73").Value = 1481.6
$ws.Range("K73").Value = 4444.799999999999
$ws.Range("M73").Value = -3508.799999999999

# Row 113
$ws.Range("H113").Value = 8112.8335
$ws.Range("I113").Value = 7293
$ws.Range("K113").Value = 7293
$ws.Range("M113").Value = -4039

# Row 125
$ws.Range("H125").Value = 1546
$ws.Range("I125").Value = 1546
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 13914
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -11454
$ws.Range("N125").ClearContents()

# Row 132
$ws.Range("H132").Value = 11020.966
$ws.Range("I132").Value = 1969.9803
$ws.Range("K132").Value = 5909.9409
$ws.Range("M132").Value = -3379.9409

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 37039910
$ws.Range("I32").Value = 38464450
$ws.Range("K32").Value = 38464450
$ws.Range("M32").Value = -38464163

# Row 45
$ws.Range("H45").Value = 2688.75
$ws.Range("I45").Value = 2252
$ws.Range("K45").Value = 2252
$ws.Range("M45").Value = -1875

# Row 102
$ws.Range("H102").Value = 8558.375
$ws.Range("J102").Value = 6000
$ws.Range("L102").Value = 6000
$ws.Range("N102").Value = -9244

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 132
$ws.Range("H132").Value = 1701.909
$ws.Range("I132").Value = 1701.909
$ws.Range("K132").Value = 5105.727000000001
$ws.Range("M132").Value = -2575.727000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2481.6316
$ws.Range("I94").Value = 2841
$ws.Range("J94").Value = 1703
$ws.Range("K94").Value = 2841
$ws.Range("L94").Value = 1703
$ws.Range("M94").Value = -2390
$ws.Range("N94").Value = -2605

# Row 134
$ws.Range("H134").Value = 1327.6552
$ws.Range("I134").Value = 1125.0714
$ws.Range("K134").Value = 3375.2142
$ws.Range("M134").Value = -840.2142000000003

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1561.5151
$ws.Range("I31").Value = 1561.5151
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1561.5151
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1266.5151
$ws.Range("N31").ClearContents()

# Row 34
$ws.Range("H34").Value = 1561.5151
$ws.Range("I34").Value = 1561.5151
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1561.5151
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1359.5151
$ws.Range("N34").ClearContents()

# Row 132
$ws.Range("H132").Value = 2002.5366
$ws.Range("I132").Value = 1841.081
$ws.Range("K132").Value = 5523.242999999999
$ws.Range("M132").Value = -2993.242999999999

# Row 134
$ws.Range("H134").Value = 2887.5334
$ws.Range("I134").Value = 2703.4167
$ws.Range("J134").Value = 3624
$ws.Range("K134").Value = 8110.250100000001
$ws.Range("L134").Value = 10872
$ws.Range("M134").Value = -5575.250100000001
$ws.Range("N134").Value = -15942

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1155.25
$ws.Range("I5").Value = 808.1539
$ws.Range("J5").Value = 1799.8572
$ws.Range("K5").Value = 2424.4617
$ws.Range("L5").Value = 5399.571599999999
$ws.Range("M5").Value = -2312.4617
$ws.Range("N5").Value = -5623.571599999999

# Row 11
$ws.Range("H11").Value = 3277.1667
$ws.Range("I11").Value = 4443.846
$ws.Range("K11").Value = 13331.538
$ws.Range("M11").Value = -13191.538

# Row 14
$ws.Range("H14").Value = 60.5
$ws.Range("I14").Value = 60.5
$ws.Range("K14").Value = 181.5
$ws.Range("M14").Value = -8.5

# Row 56
$ws.Range("H56").Value = 7807.8
$ws.Range("I56").Value = 7807.8
$ws.Range("K56").Value = 7807.8
$ws.Range("M56").Value = -7277.8

# Row 135
$ws.Range("H135").Value = 1155.25
$ws.Range("I135").Value = 808.1539
$ws.Range("J135").Value = 1799.8572
$ws.Range("K135").Value = 7273.3851
$ws.Range("L135").Value = 16198.7148
$ws.Range("M135").Value = -4738.3851
$ws.Range("N135").Value = -21268.7148

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 100609.414
$ws.Range("I80").Value = 189453
$ws.Range("J80").Value = 11765.833
$ws.Range("K80").Value = 189453
$ws.Range("L80").Value = 11765.833
$ws.Range("M80").Value = -188455
$ws.Range("N80").Value = -13761.833

# Row 83
$ws.Range("H83").Value = 100609.414
$ws.Range("I83").Value = 189453
$ws.Range("J83").Value = 11765.833
$ws.Range("K83").Value = 947265
$ws.Range("L83").Value = 58829.165
$ws.Range("M83").Value = -942273
$ws.Range("N83").Value = -68813.16500000001

# Row 102
$ws.Range("H102").Value = 3206.7856
$ws.Range("I102").Value = 2206.1765
$ws.Range("K102").Value = 2206.1765
$ws.Range("M102").Value = -584.1765

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 754.619
$ws.Range("I16").Value = 514.86664
$ws.Range("K16").Value = 514.86664
$ws.Range("M16").Value = -344.86664

# Row 88
$ws.Range("H88").Value = 31747.334
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 36096.8
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 36096.8
$ws.Range("M88").Value = -9572
$ws.Range("N88").Value = -36952.8

# Row 91
$ws.Range("H91").Value = 31747.334
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 36096.8
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 36096.8
$ws.Range("M91").Value = -8518
$ws.Range("N91").Value = -39060.8

# Row 93
$ws.Range("H93").Value = 5551.0835
$ws.Range("I93").Value = 4971.467
$ws.Range("J93").Value = 6517.1113
$ws.Range("K93").Value = 4971.467
$ws.Range("L93").Value = 6517.1113
$ws.Range("M93").Value = -3723.467
$ws.Range("N93").Value = -9013.1113

# Row 130
$ws.Range("H130").Value = 68665
$ws.Range("J130").Value = 68665
$ws.Range("L130").Value = 68665
$ws.Range("N130").Value = -78705

# Row 132
$ws.Range("H132").Value = 2808.8333
$ws.Range("I132").Value = 2501.762
$ws.Range("K132").Value = 7505.286
$ws.Range("M132").Value = -4975.286

# Row 136
$ws.Range("H136").Value = 2910.2334
$ws.Range("I136").Value = 2724.32
$ws.Range("K136").Value = 8172.960000000001
$ws.Range("M136").Value = -5622.960000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3162.9092
$ws.Range("I122").Value = 2561.5
$ws.Range("J122").Value = 4766.6665
$ws.Range("K122").Value = 7684.5
$ws.Range("L122").Value = 14299.9995
$ws.Range("M122").Value = -5234.5
$ws.Range("N122").Value = -19199.9995

# Row 126
$ws.Range("H126").Value = 1263.3478
$ws.Range("I126").Value = 1158.6875
$ws.Range("K126").Value = 3476.0625
$ws.Range("M126").Value = -1006.0625

# Row 132
$ws.Range("H132").Value = 1054.4642
$ws.Range("I132").Value = 1054.4642
$ws.Range("K132").Value = 3163.3926
$ws.Range("M132").Value = -633.3925999999997
